$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.038.30'
$ws.Range("E2").Value = '  -2.03%  '

$ws.Range("D3").Value = '1.894.79'
$ws.Range("E3").Value = '  -0.96%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = "'313.01"
$ws.Range("E5").Value = '  -0.24%  '

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").Value = "'0.5010"
$ws.Range("E7").Value = '  -0.70%  '

$ws.Range("D8").Value = "'0.3883"
$ws.Range("E8").Value = '  -2.07%  '

$ws.Range("D9").Value = "'0.09161"
$ws.Range("E9").Value = '  -4.92%  '

$ws.Range("D10").Value = "'1.125"
$ws.Range("E10").Value = '  -3.09%  '

$ws.Range("D11").Value = "'41.75"
$ws.Range("E11").Value = '  +0.29%  '

$ws.Range("D12").Value = "'6.383"
$ws.Range("E12").Value = '  -2.83%  '

$ws.Range("D13").Value = "'20.77"
$ws.Range("E13").Value = '  -1.81%  '

$ws.Range("D14").Value = '1.896.19'
$ws.Range("E14").Value = '  -0.91%  '

$ws.Range("D15").Value = "'7.283"
$ws.Range("E15").Value = '  -3.62%  '

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = "'92.26"
$ws.Range("E17").Value = '  -1.77%  '

$ws.Range("D18").Value = "'0.00001105"
$ws.Range("E18").Value = '  -2.80%  '

$ws.Range("D19").Value = "'0.06651"
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Value = "'17.88"
$ws.Range("E20").Value = '  -1.05%  '

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").Value = "'6.206"
$ws.Range("E22").Value = '  -1.17%  '

$ws.Range("D23").Value = '28.097.50'
$ws.Range("E23").Value = '  -1.98%  '

$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").Value = "'2.318"
$ws.Range("E25").Value = '  +1.35%  '

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.116.42'
$ws.Range("E26").Value = '  -1.12%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = "'2.571"
$ws.Range("E27").Value = '  -7.07%  '

$ws.Range("D28").Value = "'20.84"
$ws.Range("E28").Value = '  -2.70%  '

$ws.Range("D29").Value = "'158.10"
$ws.Range("E29").Value = '  -0.80%  '

$ws.Range("D30").Value = "'126.60"
$ws.Range("E30").Value = '  -1.84%  '

$ws.Range("D31").Value = "'1.087"
$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("D32").Value = "'0.1058"
$ws.Range("E32").Value = '  -1.74%  '

$ws.Range("D33").Value = "'5.594"
$ws.Range("E33").Value = '  -2.32%  '

$ws.Range("E34").Value = '  -0.43%  '

$ws.Range("D35").Value = "'9.558"
$ws.Range("E35").Value = '  -2.70%  '

$ws.Range("D36").Value = "'0.06589"
$ws.Range("E36").Value = '  -2.91%  '

$ws.Range("D37").Value = "'0.02399"
$ws.Range("E37").Value = '  -1.80%  '

$ws.Range("D38").Value = "'0.2200"
$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("D39").Value = "'1.222"
$ws.Range("E39").Value = '  -4.96%  '

$ws.Range("D40").Value = "'1.273"
$ws.Range("E40").Value = '  +6.19%  '

$ws.Range("D41").Value = "'0.6480"
$ws.Range("E41").Value = '  +0.86%  '

$ws.Range("D42").Value = "'4.968"
$ws.Range("E42").Value = '  -2.80%  '

$ws.Range("D43").Value = "'11.36"
$ws.Range("E43").Value = '  -2.45%  '

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("D45").Value = "'0.6075"
$ws.Range("E45").Value = '  -0.54%  '

$ws.Range("D46").Value = "'13.29"
$ws.Range("E46").Value = '  -3.70%  '

$ws.Range("D47").Value = "'1.311"
$ws.Range("E47").Value = '  +1.80%  '

$ws.Range("D48").Value = "'3.681"
$ws.Range("E48").Value = '  +0.63%  '

$ws.Range("D49").Value = "'1.995"
$ws.Range("E49").Value = '  -2.37%  '

$ws.Range("D50").Value = "'121.81"
$ws.Range("E50").Value = '  -2.70%  '

$ws.Range("D51").Value = "'1.182"
$ws.Range("E51").Value = '  -2.33%  '
